$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (lambda_ = 0.02)
$ws.Range("B2").Value = 0.09757222222222223
$ws.Range("C2").Value = 0.0003127865507069558
$ws.Range("E2").Value = 0.09256073464898067
$ws.Range("F2").Value = 0.000278872105027218
$ws.Range("G2").Value = 1080026
$ws.Range("H2").Value = 0.04452022451311358
$ws.Range("I2").Value = 0.0001984596938272205
$ws.Range("J2").Value = 1080026
$ws.Range("K2").Value = 0.01128633574592749
$ws.Range("L2").Value = 0.00008266279737757257

# Row 3 (lambda_ = 0.1)
$ws.Range("B3").Value = 0.8376977777777778
$ws.Range("C3").Value = 0.0003886732573539983
$ws.Range("E3").Value = 0.3889321431328914
$ws.Range("F3").Value = 0.0002097728379262024
$ws.Range("G3").Value = 5400883
$ws.Range("H3").Value = 0.09611243198565864
$ws.Range("I3").Value = 0.000126827842691839
$ws.Range("J3").Value = 5400883
$ws.Range("K3").Value = 0.05807639207557461
$ws.Range("L3").Value = 0.00008835601281472285

# Row 4 (lambda_ = 0.2)
$ws.Range("B4").Value = 0.9947233333333333
$ws.Range("C4").Value = 0.00007636770300152048
$ws.Range("E4").Value = 0.6260722349197325
$ws.Range("F4").Value = 0.0001472261831356737
$ws.Range("G4").Value = 10800455
$ws.Range("H4").Value = 0.1859705910538028
$ws.Range("I4").Value = 0.0001183916729308955
$ws.Range("J4").Value = 10800455
$ws.Range("K4").Value = 0.111800217070816
$ws.Range("L4").Value = 0.00009214605088138318

# Row 5 (lambda_ = 0.5)
$ws.Range("E5").Value = 0.9040107740053344
$ws.Range("F5").Value = 0.00005669147628739128
$ws.Range("G5").Value = 26999801
$ws.Range("H5").Value = 0.4556105802409433
$ws.Range("I5").Value = 0.00009584544029742404
$ws.Range("J5").Value = 26999801
$ws.Range("K5").Value = 0.2685176288677825
$ws.Range("L5").Value = 0.00009637905349392552

# Row 6 (lambda_ = 1)
$ws.Range("E6").Value = 0.9683777037037037
$ws.Range("F6").Value = 0.00002381344417715802
$ws.Range("H6").Value = 0.6077662777777778
$ws.Range("I6").Value = 0.0000664421816839611
$ws.Range("K6").Value = 0.5904021657362619
$ws.Range("L6").Value = 0.00002137981728407834

$wb.Save()
